# Update Industries (column H) flag from 1 to 0 for rows 24 through 66
# on the active worksheet, reflecting the policy data correction.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H24:H66").Value = 0
